$d = $word.ActiveDocument

# --- Step 1: locate the "git commit -m ..." paragraph that carries the
# boilerplate <w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>
# and the split " " + """ runs that need merging into a single run.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like 'git commit -m " update grouped full index "*') {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $rng = $target.Range
    $xml = $rng.WordOpenXML

    # Remove the now-superfluous paragraph-mark formatting (<w:pPr>...</w:pPr>)
    $xml = $xml.Replace('<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>', '')

    # Merge the trailing " " run and the closing quote run into one run.
    $xml = $xml.Replace('<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>"</w:t></w:r>', '<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t xml:space="preserve"> "</w:t></w:r>')

    $rng.InsertXML($xml)
}

# --- Step 2: insert two blank paragraphs right after the "git push"
# paragraph that follows the "git pull --no-rebase" line (i.e. the final
# "git push" of the second git-workflow block, right before the document's
# trailing blank paragraph).
$pushPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq "git pull --no-rebase`r") {
        $pushPara = $d.Paragraphs.Item($i + 1)
        break
    }
}

if ($pushPara -ne $null) {
    $insertPoint = $d.Range($pushPara.Range.End, $pushPara.Range.End)
    $fragment = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $insertPoint.InsertXML($fragment)
}
